$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "'26.037.66"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -2.40%  '
$ws.Range('D3').Value = "'1.667.39"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -1.73%  '
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').Value = "'216.79"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.43%  '
$ws.Range('D6').Value = "'0.5103"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.49%  '
$ws.Range('E7').Value = '  -0.06%  '
$ws.Range('E8').Value = '  +0.35%  '
$ws.Range('D9').Value = "'0.06404"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +1.68%  '
$ws.Range('D10').Value = "'21.80"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -1.69%  '
$ws.Range('D11').Value = "'0.07430"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +1.16%  '
$ws.Range('D12').Value = "'1.691.06"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.70%  '
$ws.Range('D13').Value = "'4.499"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -0.65%  '
$ws.Range('D14').Value = "'0.5842"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +0.62%  '
$ws.Range('D15').Value = "'0.000008552"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +1.13%  '
$ws.Range('D16').Value = "'64.41"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -1.77%  '
$ws.Range('D17').Value = "'26.075.18"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E18').Value = '  -1.27%  '
$ws.Range('E19').Value = '  -0.01%  '
$ws.Range('D20').Value = "'10.78"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -2.12%  '
$ws.Range('D21').Value = "'192.43"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +3.02%  '
$ws.Range('D22').Value = "'6.226"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.44%  '
$ws.Range('E23').Value = '  -0.06%  '
$ws.Range('D24').Value = "'145.15"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.34%  '
$ws.Range('E25').Value = '  +1.33%  '
$ws.Range('D26').Value = "'0.1198"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +3.42%  '
$ws.Range('D27').Value = "'15.70"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.27%  '
$ws.Range('D28').Value = "'0.06440"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +13.93%  '
$ws.Range('D29').Value = "'1.335"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -1.00%  '
$ws.Range('E30').Value = '  -1.83%  '
$ws.Range('D31').Value = "'3.543"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.86%  '
$ws.Range('D32').Value = "'3.527"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.80%  '
$ws.Range('D33').Value = "'1.647"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.11%  '
$ws.Range('E34').Value = '  -0.42%  '
$ws.Range('D35').Value = "'0.6098"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +1.44%  '
$ws.Range('D36').Value = "'2.368"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.39%  '
$ws.Range('D37').Value = "'2.705"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.13%  '
$ws.Range('D38').Value = "'6.259"
$ws.Range('D38').Style = 'Normal'
$ws.Range('D39').Value = "'0.01604"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -0.58%  '
$ws.Range('D40').Value = "'1.084.93"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -1.52%  '
$ws.Range('D41').Value = "'0.8599"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.01%  '
$ws.Range('E42').Value = '  +0.57%  '
$ws.Range('D43').Value = "'100.19"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.41%  '
$ws.Range('D44').Value = "'1.816.28"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -2.19%  '
$ws.Range('D45').Value = "'0.00000000112"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +1.72%  '
$ws.Range('D46').Value = "'56.36"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -0.68%  '
$ws.Range('D47').Value = "'1.007"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.31%  '
$ws.Range('D48').Value = "'8.054"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -1.38%  '
$ws.Range('D49').Value = "'0.05236"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.09%  '
$ws.Range('D50').Value = "'0.4286"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.78%  '
$ws.Range('D51').Value = "'6.011"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +3.91%  '
